$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replicates Python's str.title(): a letter is upper-cased if it is not
# preceded by another letter, lower-cased otherwise; non-letters pass through
# unchanged (this matters for tokens like "-DTO." -> "-Dto.").
function PyTitle($s) {
    $resultStr = ""
    $prevIsLetter = $false
    for ($k = 0; $k -lt $s.Length; $k++) {
        $ch = $s.Substring($k, 1)
        if ($ch -match '\p{L}') {
            if ($prevIsLetter) {
                $resultStr += $ch.ToLower()
            } else {
                $resultStr += $ch.ToUpper()
            }
            $prevIsLetter = $true
        } else {
            $resultStr += $ch
            $prevIsLetter = $false
        }
    }
    return $resultStr
}

# 1. Rename the header row to the short machine-friendly column names.
$ws.Cells.Item(1, 1).Value = "mx_state"
$ws.Cells.Item(1, 2).Value = "mx_municipality"
$ws.Cells.Item(1, 3).Value = "n_matriculas"
$ws.Cells.Item(1, 4).Value = "pct_matriculas"

# 2. Re-case the state (col A) and municipality (col B) text for data rows
#    2..140 from ALL CAPS to Title Case.
for ($row = 2; $row -le 140; $row++) {
    $aText = $ws.Cells.Item($row, 1).Text
    if ($aText -ne "") {
        $ws.Cells.Item($row, 1).Value = PyTitle($aText)
    }
    $bText = $ws.Cells.Item($row, 2).Text
    if ($bText -ne "") {
        $ws.Cells.Item($row, 2).Value = PyTitle($bText)
    }
}

# 3. Drop the trailing metadata/footer rows (142-146); this also shrinks the
#    sheet's used range / dimension down to A1:D140.
$ws.Range("A142:D146").EntireRow.Delete()
